# Auto-generated: apply crypto price/volume update from commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.057.93'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '3.135.16'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.97'
$ws.Range("E5").Value = '  +2.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.98'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("E8").Value = '  +10.82%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  +2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.422'
$ws.Range("E11").Value = '  +5.36%  '
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").Value = '3.676.11'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.93'
$ws.Range("E14").Value = '  +2.69%  '
$ws.Range("E15").Value = '  +5.12%  '
$ws.Range("D16").Value = '58.153.76'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("D17").Value = '3.138.76'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.22'
$ws.Range("E18").Value = '  +6.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.96'
$ws.Range("E19").Value = '  +4.14%  '
$ws.Range("E20").Value = '  +4.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.94'
$ws.Range("E21").Value = '  +6.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.33'
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("E26").Value = '  +1.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '0.0₃0880'
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.93'
$ws.Range("E29").Value = '  +10.20%  '
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("E31").Value = '  +5.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.73'
$ws.Range("E32").Value = '  +3.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.16'
$ws.Range("E33").Value = '  +6.68%  '
$ws.Range("E34").Value = '  +3.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.27'
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("E36").Value = '  +4.26%  '
$ws.Range("E37").Value = '  +10.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.48'
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.68'
$ws.Range("E39").Value = '  +6.73%  '
$ws.Range("D40").Value = '2.641.52'
$ws.Range("E40").Value = '  +9.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0679'
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("E42").Value = '  +4.95%  '
$ws.Range("E43").Value = '  +4.96%  '
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("E45").Value = '  +4.63%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("E47").Value = '  +11.80%  '
$ws.Range("E48").Value = '  +4.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.976'
$ws.Range("E49").Value = '  +3.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.23'
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.748'
$ws.Range("E51").Value = '  -0.37%  '
